$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing D:K data to F:M
$ws.Columns("D:E").Insert()

# Copy number formats/styles from column F (the old column D, now shifted) into the new D:E columns
$ws.Columns("F:F").Copy()
$ws.Columns("D:E").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new columns with the new quarter data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 200600
$ws.Range("E8").Value = 181800
$ws.Range("D9").Value = 34000
$ws.Range("E9").Value = 29800
$ws.Range("D10").Value = 166600
$ws.Range("E10").Value = 152000
$ws.Range("D12").Value = 24000
$ws.Range("E12").Value = 22700
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 138200
$ws.Range("E17").Value = 131500
$ws.Range("D18").Value = 62400
$ws.Range("E18").Value = 50300
$ws.Range("D20").Value = 2100
$ws.Range("E20").Value = 1500
$ws.Range("D21").Value = 68100
$ws.Range("E21").Value = 55200
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 64500
$ws.Range("E23").Value = 51800
$ws.Range("D24").Value = 19600
$ws.Range("E24").Value = 1700
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 44900
$ws.Range("E26").Value = 50100
$ws.Range("D27").Value = 44900
$ws.Range("E27").Value = 50100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -2100
$ws.Range("E32").Value = -1500
$ws.Range("D33").Value = 44900
$ws.Range("E33").Value = 50100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 44900
$ws.Range("E35").Value = 50100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 84200
$ws.Range("E41").Value = 76400
$ws.Range("D42").Value = 374000
$ws.Range("E42").Value = 334100
$ws.Range("D43").Value = 88200
$ws.Range("E43").Value = 74200
$ws.Range("D44").Value = 71800
$ws.Range("E44").Value = 63700
$ws.Range("D45").Value = 10500
$ws.Range("E45").Value = 13300
$ws.Range("D46").Value = 628700
$ws.Range("E46").Value = 561600
$ws.Range("D48").Value = 138500
$ws.Range("E48").Value = 134100
$ws.Range("D49").Value = 51200
$ws.Range("E49").Value = 51200
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 957100
$ws.Range("E54").Value = 887200
$ws.Range("D57").Value = 26200
$ws.Range("E57").Value = 21900
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("D59").Value = 66800
$ws.Range("E59").Value = 60100
$ws.Range("D60").Value = 93000
$ws.Range("E60").Value = 82000
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 12600
$ws.Range("E62").Value = 12200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 105600
$ws.Range("E66").Value = 94200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 325500
$ws.Range("E72").Value = 280700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 851500
$ws.Range("E76").Value = 793000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 44900
$ws.Range("E81").Value = 50100
$ws.Range("D83").Value = 3600
$ws.Range("E83").Value = 3400
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 74600
$ws.Range("E89").Value = 60600
$ws.Range("D91").Value = -8100
$ws.Range("E91").Value = -12200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -65800
$ws.Range("E94").Value = -49700
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -1100
$ws.Range("E100").Value = 4100
$ws.Range("D101").Value = 100
$ws.Range("E101").Value = 100
$ws.Range("D102").Value = 7800
$ws.Range("E102").Value = 15100

# Rows 47 and 52 received restated figures beyond the simple column shift; set the full row explicitly
$ws.Range("D47").Value = 40800
$ws.Range("E47").Value = 22300
$ws.Range("F47").Value = 20700
$ws.Range("G47").Value = 50100
$ws.Range("H47").Value = 62100
$ws.Range("I47").Value = 68600
$ws.Range("J47").Value = 45300
$ws.Range("K47").Value = 47100
$ws.Range("L47").Value = 18200
$ws.Range("M47").Value = 4000
$ws.Range("D52").Value = 98000
$ws.Range("E52").Value = 118000
$ws.Range("F52").Value = 116900
$ws.Range("G52").Value = 72300
$ws.Range("H52").Value = 76300
$ws.Range("I52").Value = 107900
$ws.Range("J52").Value = 114500
$ws.Range("K52").Value = 43000
$ws.Range("L52").Value = 43600
$ws.Range("M52").Value = 53000
